# Add two new datapoint rows (U_DC_NE, U_DC_PE) right after the existing
# "U_DC[1..x]" row (row 42) in the Worksheet sheet, pushing the following
# rows (STATE[1..x], ERROR[1..x], QS_TX, QS_RX) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 43/44 (shifts old rows 43-46 down to 45-48)
$ws.Range("A43:A44").EntireRow.Insert()

# Row 43: U_DC_NE
$ws.Range("A43").Value = "datapoints"
$ws.Range("B43").Value = "U_DC_NE"
$ws.Range("C43").Value = "V"
$ws.Range("D43").Value = "Voltage DC negative pole to earth"

# Row 44: U_DC_PE
$ws.Range("A44").Value = "datapoints"
$ws.Range("B44").Value = "U_DC_PE"
$ws.Range("C44").Value = "V"
$ws.Range("D44").Value = "Voltage DC positive pole to earth"

# The two new description cells carry an explicit (visually-default) font,
# matching the distinct cell style used for these rows in the source file.
$ws.Range("D43:D44").Font.Bold = $false

$null = $ws.Range("D45").Select()
